$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.080499649047852
$ws.Range("B1").Value = 2.387287139892578
$ws.Range("C1").Value = 2.045372247695923
$ws.Range("D1").Value = 2.162761688232422
$ws.Range("E1").Value = 2.568912506103516
